# "timer and screen setup finished"
#
# A new "System Timer" entry (8 B) is carved out of the former "Empty"
# block, so the old blank spacer row above "Free Memory" is removed to
# keep the table the same overall length. The "Free Memory" size note
# also changes from an exact "3068 MB" to "<=3068 MB".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new "System Timer" row by inserting a blank row
# right after the "GUI Sheet Controller" row (pushes "Empty" ... "Stack"
# and the spacer row down by one).
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).RowHeight = 35.25

# Fill in the new System Timer row.
$ws.Range("A5").Value = "0x0026a414 - 0x0026a41b"
$ws.Range("B5").Value = "System Timer"
$ws.Range("C5").Value = "8 B"

# The "Empty" block now starts right after the timer and shrank by 8 B.
$ws.Range("A6").Value = "0x0026a41c - 0x0026f7ff"
$ws.Range("C6").Value = "21476 B"

# Drop the now-duplicated blank spacer row (originally row 10, shifted to
# row 11 by the insert above) so "Free Memory" lands back on row 11.
$ws.Rows.Item(11).Delete()

# "Free Memory" size note is now an upper bound rather than an exact value.
$ws.Range("C11").Value = "<=3068 MB"

# Reset the selection to C7 instead of A11 (the view already re-scrolls
# to the top on save, matching the dropped topLeftCell="A4").
$ws.Range("C7").Select()
